$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

$startSerial = 44867
for ($i = 0; $i -lt 19; $i++) {
    $row = 22 + $i
    $prevRow = $row - 1
    # Pre-touch the new row with a throwaway write so the engine registers
    # it in the dependency graph / used-range *before* Copy() runs - doing
    # Copy() first on a row beyond the current dimension leaves it outside
    # dependent SUM() ranges.
    $ws2.Range("A$row").Value = 0
    $ws2.Range("A${prevRow}:J${prevRow}").Copy($ws2.Range("A${row}:J${row}"))

    $ws2.Range("A$row").Value = $startSerial + $i
    $ws2.Range("B$row").Value = 20
    $ws2.Range("C$row").Value = 11
    $ws2.Range("D$row").Value = 9
    $ws2.Range("E$row").Formula = "=D$row/B$row*100"
    $ws2.Range("F$row").Value = 1
    $ws2.Range("G$row").Value = 3
    $ws2.Range("H$row").Value = 3
    $ws2.Range("I$row").Value = 1
    $ws2.Range("J$row").Value = 1
}

$ws2.Range("A25").Select()

Write-Output "done"
